$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Chocolate Treasures"
$ws.Range("A6").Value = "Water Seal Key Treasures"
$ws.Range("A7").Value = "Kachi Katchin Treasures"

$ws.Range("A4").Select()
